# Apply updated cryptocurrency price / 1h-volume-change figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.224.85'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '3.804.53'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '701.79'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.15'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("D7").Value = '3.803.64'
$ws.Range("E7").Value = '  -1.09%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.47'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.87'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").Value = '4.447.66'
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("D16").Value = '3.802.62'
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").Value = '71.277.73'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.16'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("E19").Value = '  -0.47%  '
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '513.50'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.52'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.713'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.99'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("E25").Value = '  -3.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.12'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").Value = '3.956.08'
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.31'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("E30").Value = '  -4.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.01'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.34'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.175'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.10'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = '3.768.12'
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  -2.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.30'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.36'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("E42").Value = '  -2.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.26'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '173.17'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000308'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.47'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '421.79'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("E50").Value = '  -1.97%  '
$ws.Range("E51").Value = '  -1.14%  '
